$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.182.81"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.51%  "
$ws.Range("D3").Value = "'1.879.59"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.42%  "
$ws.Range("D4").Value = "'1.003"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").Value = "'313.00"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.08%  "
$ws.Range("E6").Value = "  -0.14%  "
$ws.Range("D7").Value = "'0.5128"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.82%  "
$ws.Range("D8").Value = "'0.3898"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.15%  "
$ws.Range("D9").Value = "'0.08369"
$ws.Range("D9").Style = "Normal"
$ws.Range("E10").Value = "  +0.42%  "
$ws.Range("D11").Value = "'41.59"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.11%  "
$ws.Range("D12").Value = "'6.228"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.07%  "
$ws.Range("D13").Value = "'20.74"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Value = "'1.875.18"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.11%  "
$ws.Range("D15").Value = "'7.299"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.29%  "
$ws.Range("D16").Value = "'1.002"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.24%  "
$ws.Range("D17").Value = "'0.00001108"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.13%  "
$ws.Range("D18").Value = "'91.36"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.51%  "
$ws.Range("D19").Value = "'0.06649"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.19%  "
$ws.Range("D20").Value = "'17.76"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.18%  "
$ws.Range("E21").Value = "  -0.04%  "
$ws.Range("D22").Value = "'6.051"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.10%  "
$ws.Range("D23").Value = "'28.222.19"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.51%  "
$ws.Range("D24").Value = "'11.20"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.43%  "
$ws.Range("D25").Value = "'2.259"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.32%  "
$ws.Range("D26").Value = "'2.091.04"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.21%  "
$ws.Range("D27").Value = "'2.513"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.08%  "
$ws.Range("D28").Value = "'158.69"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.97%  "
$ws.Range("E29").Value = "  +0.24%  "
$ws.Range("D30").Value = "'125.42"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.28%  "
$ws.Range("E31").Value = "  +0.94%  "
$ws.Range("E32").Value = "  -0.27%  "
$ws.Range("D33").Value = "'5.895"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +5.21%  "
$ws.Range("D34").Value = "'3.583"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.62%  "
$ws.Range("D35").Value = "'9.719"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.38%  "
$ws.Range("D36").Value = "'0.02458"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.28%  "
$ws.Range("D37").Value = "'0.06552"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.58%  "
$ws.Range("D38").Value = "'0.2196"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.38%  "
$ws.Range("D39").Value = "'1.212"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.17%  "
$ws.Range("D40").Value = "'0.6522"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.64%  "
$ws.Range("D41").Value = "'5.022"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.84%  "
$ws.Range("E42").Value = "  -1.35%  "
$ws.Range("D43").Value = "'11.30"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.18%  "
$ws.Range("D44").Value = "'0.6119"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.13%  "
$ws.Range("D45").Value = "'13.08"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.33%  "
$ws.Range("D46").Value = "'1.291"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.41%  "
$ws.Range("D47").Value = "'3.676"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.29%  "
$ws.Range("D48").Value = "'2.018"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.36%  "
$ws.Range("E49").Value = "  +1.02%  "
$ws.Range("D50").Value = "'121.75"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.74%  "
$ws.Range("D51").Value = "'78.22"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.71%  "
